$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A52").Value = 43953
$ws.Range("B52").Value = 36318
$ws.Range("C52").Value = 1983
$ws.Range("D52").Value = 78
$ws.Range("E52").Value = 4451

$table = $ws.ListObjects.Item("Table3")
$table.Resize($ws.Range("A1:E52"))

$ws.Range("E52").Select()
